$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) New column J ("Meta") header, matching existing header look (bold,
#    centered, bordered) copied from A1, with explicit (no-op) fill touch so
#    the resulting style carries applyFill like the real edit did.
# ---------------------------------------------------------------------------
$ws.Range("J1").Value = "Meta"
$ws.Range("A1").Copy()
$ws.Range("J1").PasteSpecial(-4122)
$ws.Range("J1").Interior.ThemeColor = 1
$ws.Range("J1").Interior.Pattern = -4142

# ---------------------------------------------------------------------------
# 2) Build the needed body styles on out-of-the-way scratch cells first, so
#    that applying them later to the real ranges (via copy/paste-format)
#    never forces the engine to fabricate extra transient style records.
# ---------------------------------------------------------------------------
$s_body       = $ws.Range("Z1")   # default font + border
$s_bodyDate   = $ws.Range("Z2")   # default font + border + short-date fmt
$s_body2      = $ws.Range("Z3")   # font2 (explicit black) + border
$s_body2Date  = $ws.Range("Z4")   # font2 (explicit black) + border + date

$s_body.Borders.LineStyle = 1

$s_bodyDate.Borders.LineStyle = 1
$s_bodyDate.NumberFormat = "mm-dd-yy"

$s_body2.Borders.LineStyle = 1
$s_body2.Font.Color = 0

$s_body2Date.Borders.LineStyle = 1
$s_body2Date.Font.Color = 0
$s_body2Date.NumberFormat = "mm-dd-yy"

# ---------------------------------------------------------------------------
# 3) Apply the "body" style to the existing data block (rows 2-12) plus the
#    brand-new column J, and the "body date" style to column B of that block.
# ---------------------------------------------------------------------------
$s_body.Copy()
$ws.Range("A2:A12").PasteSpecial(-4122)
$ws.Range("C2:J12").PasteSpecial(-4122)

$s_bodyDate.Copy()
$ws.Range("B2:B12").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 4) Populate the new rows (13 and 14) with the pasted-in measurements.
# ---------------------------------------------------------------------------
$ws.Range("A13").Value = "Antoniazzi"
$ws.Range("B13").Value = 46055
$ws.Range("C13").Value = 71.2
$ws.Range("D13").Value = 23.5
$ws.Range("E13").Value = 22
$ws.Range("F13").Value = 36.8
$ws.Range("G13").Value = 1613
$ws.Range("H13").Value = 40
$ws.Range("I13").Value = 7

$ws.Range("A14").Value = "Isabely"
$ws.Range("B14").Value = 46144
$ws.Range("C14").Value = 88.1
$ws.Range("D14").Value = 30.4
$ws.Range("E14").Value = 47.2
$ws.Range("F14").Value = 22.6
$ws.Range("G14").Value = 1606
$ws.Range("H14").Value = 54
$ws.Range("I14").Value = 8

$s_body2.Copy()
$ws.Range("A13:A14").PasteSpecial(-4122)
$ws.Range("C13:J14").PasteSpecial(-4122)

$s_body2Date.Copy()
$ws.Range("B13:B14").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 5) The new "Meta" column only has a value for the first data row.
# ---------------------------------------------------------------------------
$ws.Range("J2").Value = 70

# ---------------------------------------------------------------------------
# 6) Clean up the scratch cells used to build the styles.
# ---------------------------------------------------------------------------
$ws.Range("Z1:Z4").Clear()

# ---------------------------------------------------------------------------
# 7) Cosmetic sheet-level tweaks seen in the real edit.
# ---------------------------------------------------------------------------
$ws.Columns("B").ColumnWidth = 25 + (1/12)
$ws.Range("N18").Select()

Write-Host "Done"
